$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 0.1481638611118683
$ws.Cells.Item(2, 4).Value = 0.09410957921863172
$ws.Cells.Item(2, 5).Value = 0.1393632720748528
$ws.Cells.Item(2, 6).Value = 2.504620371392278
$ws.Cells.Item(2, 7).Value = 1.799140899952462
$ws.Cells.Item(2, 8).Value = 1.567977822701863
$ws.Cells.Item(2, 9).Value = 1.785663279673443
$ws.Cells.Item(2, 10).Value = 0.2125269747954945
$ws.Cells.Item(2, 11).Value = 2.310733107587055
$ws.Cells.Item(2, 12).Value = 0.2141147150665432
$ws.Cells.Item(2, 14).Value = 1.584938913018448
$ws.Cells.Item(3, 3).Value = 0.1458957587269794
$ws.Cells.Item(3, 4).Value = 0.09232494463840624
$ws.Cells.Item(3, 5).Value = 0.1385824311786337
$ws.Cells.Item(3, 6).Value = 2.508053576015726
$ws.Cells.Item(3, 7).Value = 1.801340132458861
$ws.Cells.Item(3, 8).Value = 1.575535957438319
$ws.Cells.Item(3, 9).Value = 1.788478908637337
$ws.Cells.Item(3, 10).Value = 0.2124992065843045
$ws.Cells.Item(3, 11).Value = 2.17386656949293
$ws.Cells.Item(3, 12).Value = 0.2137970777133376
$ws.Cells.Item(3, 14).Value = 1.599329280835818
$ws.Cells.Item(4, 3).Value = 0.1445518762630513
$ws.Cells.Item(4, 4).Value = 0.09125138132426969
$ws.Cells.Item(4, 5).Value = 0.1381524701775554
$ws.Cells.Item(4, 6).Value = 2.511525539358104
$ws.Cells.Item(4, 7).Value = 1.803816231146342
$ws.Cells.Item(4, 8).Value = 1.58093262398792
$ws.Cells.Item(4, 9).Value = 1.791192849312914
$ws.Cells.Item(4, 10).Value = 0.2125692143848781
$ws.Cells.Item(4, 11).Value = 2.090560682656587
$ws.Cells.Item(4, 12).Value = 0.2136800019165932
$ws.Cells.Item(4, 14).Value = 1.608766439860702
$ws.Cells.Item(5, 3).Value = 0.1440165468228969
$ws.Cells.Item(5, 4).Value = 0.09081952919055425
$ws.Cells.Item(5, 5).Value = 0.137989730390057
$ws.Cells.Item(5, 6).Value = 2.513283108557602
$ws.Cells.Item(5, 7).Value = 1.805107875645874
$ws.Cells.Item(5, 8).Value = 1.583321824385976
$ws.Cells.Item(5, 9).Value = 1.792546307337368
$ws.Cells.Item(5, 10).Value = 0.2126196528233777
$ws.Cells.Item(5, 11).Value = 2.056798062195867
$ws.Cells.Item(5, 12).Value = 0.2136519279173967
$ws.Cells.Item(5, 14).Value = 1.612763495069309
$ws.Cells.Item(6, 3).Value = 0.143928401175188
$ws.Cells.Item(6, 4).Value = 0.09074816217341919
$ws.Cells.Item(6, 5).Value = 0.1379634618182983
$ws.Cells.Item(6, 6).Value = 2.513595643195785
$ws.Cells.Item(6, 7).Value = 1.8053394078336
$ws.Cells.Item(6, 8).Value = 1.583730023663179
$ws.Cells.Item(6, 9).Value = 1.79278599051873
$ws.Cells.Item(6, 10).Value = 0.2126293519277667
$ws.Cells.Item(6, 11).Value = 2.0512030340401
$ws.Cells.Item(6, 12).Value = 0.2136484531541925
$ws.Cells.Item(6, 14).Value = 1.613436346619892
$ws.Cells.Item(7, 3).Value = 0.1445446066887257
$ws.Cells.Item(7, 4).Value = 0.09124553433974825
$ws.Cells.Item(7, 5).Value = 0.138150224868145
$ws.Cells.Item(7, 6).Value = 2.511547855263046
$ws.Cells.Item(7, 7).Value = 1.803832507046693
$ws.Cells.Item(7, 8).Value = 1.580964076319262
$ws.Cells.Item(7, 9).Value = 1.791210100671563
$ws.Cells.Item(7, 10).Value = 0.2125698058796885
$ws.Cells.Item(7, 11).Value = 2.090104596322874
$ws.Cells.Item(7, 12).Value = 0.2136795437538694
$ws.Cells.Item(7, 14).Value = 1.608819732651391
$ws.Cells.Item(8, 3).Value = 0.1473717362805331
$ws.Cells.Item(8, 4).Value = 0.09348965617248695
$ws.Cells.Item(8, 5).Value = 0.139083786924342
$ws.Cells.Item(8, 6).Value = 2.505520884517125
$ws.Cells.Item(8, 7).Value = 1.799665256687859
$ws.Cells.Item(8, 8).Value = 1.57042692585641
$ws.Cells.Item(8, 9).Value = 1.786429510072125
$ws.Cells.Item(8, 10).Value = 0.2124993419410401
$ws.Cells.Item(8, 11).Value = 2.263390689267794
$ws.Cells.Item(8, 12).Value = 0.2139890333429904
$ws.Cells.Item(8, 14).Value = 1.58977595028778
$ws.Cells.Item(9, 3).Value = 0.1533003704920191
$ws.Cells.Item(9, 4).Value = 0.09806472885688322
$ws.Cells.Item(9, 5).Value = 0.1413059487735673
$ws.Cells.Item(9, 6).Value = 2.504539314294831
$ws.Cells.Item(9, 7).Value = 1.800449138734294
$ws.Cells.Item(9, 8).Value = 1.555765956234453
$ws.Cells.Item(9, 9).Value = 1.78488315798343
$ws.Cells.Item(9, 10).Value = 0.2130513522610045
$ws.Cells.Item(9, 11).Value = 2.608957320897218
$ws.Cells.Item(9, 12).Value = 0.2152132706200618
$ws.Cells.Item(9, 14).Value = 1.557198541196563
$ws.Cells.Item(10, 3).Value = 0.157888405100195
$ws.Cells.Item(10, 4).Value = 0.1015302327964918
$ws.Cells.Item(10, 5).Value = 0.1431760084668099
$ws.Cells.Item(10, 6).Value = 2.510449797784929
$ws.Cells.Item(10, 7).Value = 1.806520622307033
$ws.Cells.Item(10, 8).Value = 1.548661165496299
$ws.Cells.Item(10, 9).Value = 1.788538795477763
$ws.Cells.Item(10, 10).Value = 0.2138772485301743
$ws.Cells.Item(10, 11).Value = 2.866320677326541
$ws.Cells.Item(10, 12).Value = 0.2164877919457879
$ws.Cells.Item(10, 14).Value = 1.536164210111103
$ws.Cells.Item(11, 3).Value = 0.1600256289665083
$ws.Cells.Item(11, 4).Value = 0.1031289923987941
$ws.Cells.Item(11, 5).Value = 0.1440780756842628
$ws.Cells.Item(11, 6).Value = 2.514584378481658
$ws.Cells.Item(11, 7).Value = 1.810484017901842
$ws.Cells.Item(11, 8).Value = 1.546226948277564
$ws.Cells.Item(11, 9).Value = 1.791246758873683
$ws.Cells.Item(11, 10).Value = 0.2143441828169514
$ws.Cells.Item(11, 11).Value = 2.984151908180593
$ws.Cells.Item(11, 12).Value = 0.2171488102335317
$ws.Cells.Item(11, 14).Value = 1.527223798356935
$ws.Cells.Item(12, 3).Value = 0.1608420973962694
$ws.Cells.Item(12, 4).Value = 0.103737562828556
$ws.Cells.Item(12, 5).Value = 0.144427026568561
$ws.Cells.Item(12, 6).Value = 2.516358356074633
$ws.Cells.Item(12, 7).Value = 1.812158219853984
$ws.Cells.Item(12, 8).Value = 1.545420025940672
$ws.Cells.Item(12, 9).Value = 1.792422780562234
$ws.Cells.Item(12, 10).Value = 0.2145341072007909
$ws.Cells.Item(12, 11).Value = 3.028879146411725
$ws.Cells.Item(12, 12).Value = 0.2174107760586992
$ws.Cells.Item(12, 14).Value = 1.52392857315909
$ws.Cells.Item(13, 3).Value = 0.160665939294347
$ws.Cells.Item(13, 4).Value = 0.1036063567702001
$ws.Cells.Item(13, 5).Value = 0.1443515469793404
$ws.Cells.Item(13, 6).Value = 2.515967027766209
$ws.Cells.Item(13, 7).Value = 1.811789931377859
$ws.Cells.Item(13, 8).Value = 1.545588700248544
$ws.Cells.Item(13, 9).Value = 1.792162801487564
$ws.Cells.Item(13, 10).Value = 0.2144926208618188
$ws.Cells.Item(13, 11).Value = 3.019241595576261
$ws.Cells.Item(13, 12).Value = 0.2173538392311016
$ws.Cells.Item(13, 14).Value = 1.524634242345762
$ws.Cells.Item(14, 3).Value = 0.1600926573795078
$ws.Cells.Item(14, 4).Value = 0.1031789969086816
$ws.Cells.Item(14, 5).Value = 0.1441066367911468
$ws.Cells.Item(14, 6).Value = 2.514726147775775
$ws.Cells.Item(14, 7).Value = 1.810618277953296
$ws.Cells.Item(14, 8).Value = 1.546158259443558
$ws.Cells.Item(14, 9).Value = 1.791340491379785
$ws.Cells.Item(14, 10).Value = 0.214359545405344
$ws.Cells.Item(14, 11).Value = 2.987829506968239
$ws.Cells.Item(14, 12).Value = 0.2171701289472381
$ws.Cells.Item(14, 14).Value = 1.526950888248898
$ws.Cells.Item(15, 3).Value = 0.1597424348741896
$ws.Cells.Item(15, 4).Value = 0.1029176360447934
$ws.Cells.Item(15, 5).Value = 0.1439575795511061
$ws.Cells.Item(15, 6).Value = 2.513993211547785
$ws.Cells.Item(15, 7).Value = 1.809923199085603
$ws.Cells.Item(15, 8).Value = 1.546522093703175
$ws.Cells.Item(15, 9).Value = 1.790856421106156
$ws.Cells.Item(15, 10).Value = 0.2142797393430484
$ws.Cells.Item(15, 11).Value = 2.968602599829694
$ws.Cells.Item(15, 12).Value = 0.2170591177595682
$ws.Cells.Item(15, 14).Value = 1.528381661901612
$ws.Cells.Item(16, 3).Value = 0.1577497358514961
$ws.Cells.Item(16, 4).Value = 0.1014261942670913
$ws.Cells.Item(16, 5).Value = 0.1431180874157079
$ws.Cells.Item(16, 6).Value = 2.510208720810141
$ws.Cells.Item(16, 7).Value = 1.806285831042061
$ws.Cells.Item(16, 8).Value = 1.548836312909799
$ws.Cells.Item(16, 9).Value = 1.788382877718888
$ws.Cells.Item(16, 10).Value = 0.2138485679986744
$ws.Cells.Item(16, 11).Value = 2.858635199386526
$ws.Cells.Item(16, 12).Value = 0.2164462249962114
$ws.Cells.Item(16, 14).Value = 1.536761126253182
$ws.Cells.Item(17, 3).Value = 0.1565400718368721
$ws.Cells.Item(17, 4).Value = 0.1005169152129497
$ws.Cells.Item(17, 5).Value = 0.1426162219476836
$ws.Cells.Item(17, 6).Value = 2.508257654483288
$ws.Cells.Item(17, 7).Value = 1.804362556720605
$ws.Cells.Item(17, 8).Value = 1.55046044404645
$ws.Cells.Item(17, 9).Value = 1.787133303623264
$ws.Cells.Item(17, 10).Value = 0.2136074156691663
$ws.Cells.Item(17, 11).Value = 2.791366144586675
$ws.Cells.Item(17, 12).Value = 0.216091020206477
$ws.Cells.Item(17, 14).Value = 1.542062538441336
$ws.Cells.Item(18, 3).Value = 0.1558490244736674
$ws.Cells.Item(18, 4).Value = 0.09999602241570926
$ws.Cells.Item(18, 5).Value = 0.1423323988904137
$ws.Cells.Item(18, 6).Value = 2.507271521555879
$ws.Cells.Item(18, 7).Value = 1.803369390921972
$ws.Cells.Item(18, 8).Value = 1.55146968415977
$ws.Cells.Item(18, 9).Value = 1.786512923039666
$ws.Cells.Item(18, 10).Value = 0.2134772989265272
$ws.Cells.Item(18, 11).Value = 2.75274594971728
$ws.Cells.Item(18, 12).Value = 0.215894363442473
$ws.Cells.Item(18, 14).Value = 1.545170901440819
$ws.Cells.Item(19, 3).Value = 0.1556158600035644
$ws.Cells.Item(19, 4).Value = 0.09982001933874329
$ws.Cells.Item(19, 5).Value = 0.1422371328737597
$ws.Cells.Item(19, 6).Value = 2.506960992153125
$ws.Cells.Item(19, 7).Value = 1.803052520925178
$ws.Cells.Item(19, 8).Value = 1.551824286807232
$ws.Cells.Item(19, 9).Value = 1.786319753823236
$ws.Cells.Item(19, 10).Value = 0.2134347189595118
$ws.Cells.Item(19, 11).Value = 2.739682087910182
$ws.Cells.Item(19, 12).Value = 0.2158290932057412
$ws.Cells.Item(19, 14).Value = 1.54623349632481
$ws.Cells.Item(20, 3).Value = 0.1566683546111847
$ws.Cells.Item(20, 4).Value = 0.1006134924914761
$ws.Cells.Item(20, 5).Value = 0.1426691460160683
$ws.Cells.Item(20, 6).Value = 2.508451263416646
$ws.Cells.Item(20, 7).Value = 1.804555588011169
$ws.Cells.Item(20, 8).Value = 1.550279780695206
$ws.Cells.Item(20, 9).Value = 1.787256142746017
$ws.Cells.Item(20, 10).Value = 0.2136321980418785
$ws.Cells.Item(20, 11).Value = 2.798519693784044
$ws.Cells.Item(20, 12).Value = 0.2161280410695881
$ws.Cells.Item(20, 14).Value = 1.541492074407955
$ws.Cells.Item(21, 3).Value = 0.1602608506258605
$ws.Cells.Item(21, 4).Value = 0.1033044376316639
$ws.Cells.Item(21, 5).Value = 0.1441783733965032
$ws.Cells.Item(21, 6).Value = 2.515084967950685
$ws.Cells.Item(21, 7).Value = 1.810957711622336
$ws.Cells.Item(21, 8).Value = 1.545987847608018
$ws.Cells.Item(21, 9).Value = 1.791577934780875
$ws.Cells.Item(21, 10).Value = 0.2143982772817452
$ws.Cells.Item(21, 11).Value = 2.997053103893734
$ws.Cells.Item(21, 12).Value = 0.2172237730893016
$ws.Cells.Item(21, 14).Value = 1.526267982205077
$ws.Cells.Item(22, 3).Value = 0.16265039324432
$ws.Cells.Item(22, 4).Value = 0.1050814920522924
$ws.Cells.Item(22, 5).Value = 0.1452076132466829
$ws.Cells.Item(22, 6).Value = 2.520634802084516
$ws.Cells.Item(22, 7).Value = 1.816152560294967
$ws.Cells.Item(22, 8).Value = 1.543852398260753
$ws.Cells.Item(22, 9).Value = 1.795280285113478
$ws.Cells.Item(22, 10).Value = 0.2149753436971906
$ws.Cells.Item(22, 11).Value = 3.127429784821004
$ws.Cells.Item(22, 12).Value = 0.2180078042998304
$ws.Cells.Item(22, 14).Value = 1.516844576879414
$ws.Cells.Item(23, 3).Value = 0.1613712580536202
$ws.Cells.Item(23, 4).Value = 0.1041313804490898
$ws.Cells.Item(23, 5).Value = 0.1446543752366125
$ws.Cells.Item(23, 6).Value = 2.51756151109015
$ws.Cells.Item(23, 7).Value = 1.813287299577865
$ws.Cells.Item(23, 8).Value = 1.544930813688325
$ws.Cells.Item(23, 9).Value = 1.79322384942931
$ws.Cells.Item(23, 10).Value = 0.2146603667974887
$ws.Cells.Item(23, 11).Value = 3.057788729821141
$ws.Cells.Item(23, 12).Value = 0.2175831479178498
$ws.Cells.Item(23, 14).Value = 1.52182586323611
$ws.Cells.Item(24, 3).Value = 0.1566103442457205
$ws.Cells.Item(24, 4).Value = 0.100569824095416
$ws.Cells.Item(24, 5).Value = 0.1426452043842801
$ws.Cells.Item(24, 6).Value = 2.508363310567248
$ws.Cells.Item(24, 7).Value = 1.804467968010925
$ws.Cells.Item(24, 8).Value = 1.550361223397601
$ws.Cells.Item(24, 9).Value = 1.787200301876084
$ws.Cells.Item(24, 10).Value = 0.21362096737613
$ws.Cells.Item(24, 11).Value = 2.795285406861694
$ws.Cells.Item(24, 12).Value = 0.2161112804017904
$ws.Cells.Item(24, 14).Value = 1.541749792698184
$ws.Cells.Item(25, 3).Value = 0.1516555706112968
$ws.Cells.Item(25, 4).Value = 0.0968086023192285
$ws.Cells.Item(25, 5).Value = 0.140663012971963
$ws.Cells.Item(25, 6).Value = 2.503642090415113
$ws.Cells.Item(25, 7).Value = 1.79927460480036
$ws.Cells.Item(25, 8).Value = 1.559088787949847
$ws.Cells.Item(25, 9).Value = 1.784461420054896
$ws.Cells.Item(25, 10).Value = 0.212828174659073
$ws.Cells.Item(25, 11).Value = 2.514859818393177
$ws.Cells.Item(25, 12).Value = 0.2152132706200618
$ws.Cells.Item(25, 14).Value = 1.565501866965739
